$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "increase/decrease" ratio columns:
#   K = NewExpr / RefExpr  (E / D)
#   L = NewFlux  / RefFlux (H / F)
for ($r = 2; $r -le 98; $r++) {
    $ws.Range("K$r").Formula = "=E$r/D$r"
    $ws.Range("L$r").Formula = "=H$r/F$r"
}

# Match the active selection left by the author when the file was saved.
[void]$ws.Range("P6").Select()
